$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.717.00'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '1.630.61'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  -0.77%  '
$ws.Range('D5').Value = '''219.05'
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = '''1.01'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = '''0.496'
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('E8').Value = '  -1.62%  '
$ws.Range('D9').Value = '''0.0619'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').Value = '''18.88'
$ws.Range('E10').Value = '  -1.54%  '
$ws.Range('D11').Value = '''0.0843'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '1.857.95'
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('D13').Value = '1.624.53'
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').Value = '''0.520'
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('D16').Value = '''64.01'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '26.701.82'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('D19').Value = '''212.48'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').Value = '  -0.67%  '
$ws.Range('D21').Value = '''4.31'
$ws.Range('E21').Value = '  -1.16%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '''6.20'
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = '''2.33'
$ws.Range('E23').Value = '  -3.53%  '
$ws.Range('D24').Value = '''8.96'
$ws.Range('E24').Value = '  -4.73%  '
$ws.Range('D25').Value = '''147.70'
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('D28').Value = '''6.99'
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('D29').Value = '''15.56'
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').Value = '''0.0497'
$ws.Range('E30').Value = '  -4.22%  '
$ws.Range('E31').Value = '  +0.71%  '
$ws.Range('D32').Value = '''3.36'
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '''1.52'
$ws.Range('E34').Value = '  -1.05%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.250.71'
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('D36').Value = '''2.44'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = '''0.0174'
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('D38').Value = '''0.522'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').Value = '''0.800'
$ws.Range('E40').Value = '  -3.80%  '
$ws.Range('D41').Value = '''0.800'
$ws.Range('E41').Value = '  -1.88%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''5.25'
$ws.Range('E42').Value = '  -2.28%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.767.88'
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '''2.14'
$ws.Range('E44').Value = '  -4.84%  '
$ws.Range('D45').Value = '''91.60'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = '''59.46'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('E47').Value = '  -3.08%  '
$ws.Range('D48').Value = '''0.0515'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₇0965'
$ws.Range('E49').Value = '  -7.43%  '
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('E51').Value = '  -1.09%  '
